# Update NATMI LR-pair edge-weight metrics for Lamb2-Rpsa with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 34.29298266666667
$ws.Range("H2").Value = 102.878948
$ws.Range("I2").Value = 0.1703377185274763
$ws.Range("J2").Value = 0.1703377185274763
$ws.Range("M2").Value = 20.56839166666667
$ws.Range("N2").Value = 61.705175
$ws.Range("O2").Value = 0.1304525281245593
$ws.Range("P2").Value = 0.1304525281245593
$ws.Range("Q2").Value = 705.3514989062111
$ws.Range("R2").Value = 6348.1634901559
$ws.Range("S2").Value = 0.02222098601687886
$ws.Range("T2").Value = 0.02222098601687886

# Row 3
$ws.Range("G3").Value = 34.29298266666667
$ws.Range("H3").Value = 102.878948
$ws.Range("I3").Value = 0.1703377185274763
$ws.Range("J3").Value = 0.1703377185274763
$ws.Range("O3").Value = 0.6526310778549473
$ws.Range("P3").Value = 0.6526310778549473
$ws.Range("Q3").Value = 3528.74961961814
$ws.Range("R3").Value = 31758.74657656326
$ws.Range("S3").Value = 0.1111676888419395
$ws.Range("T3").Value = 0.1111676888419395

# Row 4
$ws.Range("G4").Value = 34.29298266666667
$ws.Range("H4").Value = 102.878948
$ws.Range("I4").Value = 0.1703377185274763
$ws.Range("J4").Value = 0.1703377185274763
$ws.Range("O4").Value = 0.2169163940204933
$ws.Range("P4").Value = 0.2169163940204934
$ws.Range("Q4").Value = 1172.858095272748
$ws.Range("R4").Value = 10555.72285745474
$ws.Range("S4").Value = 0.03694904366865793
$ws.Range("T4").Value = 0.03694904366865794

# Row 5
$ws.Range("I5").Value = 0.5338339874103145
$ws.Range("J5").Value = 0.5338339874103145
$ws.Range("M5").Value = 20.56839166666667
$ws.Range("N5").Value = 61.705175
$ws.Range("O5").Value = 0.1304525281245593
$ws.Range("P5").Value = 0.1304525281245593
$ws.Range("Q5").Value = 2210.553284627955
$ws.Range("R5").Value = 19894.9795616516
$ws.Range("S5").Value = 0.06963999325648969
$ws.Range("T5").Value = 0.0696399932564897

# Row 6
$ws.Range("I6").Value = 0.5338339874103145
$ws.Range("J6").Value = 0.5338339874103145
$ws.Range("O6").Value = 0.6526310778549473
$ws.Range("P6").Value = 0.6526310778549473
$ws.Range("S6").Value = 0.3483966505991979
$ws.Range("T6").Value = 0.3483966505991979

# Row 7
$ws.Range("I7").Value = 0.5338339874103145
$ws.Range("J7").Value = 0.5338339874103145
$ws.Range("O7").Value = 0.2169163940204933
$ws.Range("P7").Value = 0.2169163940204934
$ws.Range("S7").Value = 0.1157973435546269
$ws.Range("T7").Value = 0.1157973435546269

# Row 8
$ws.Range("G8").Value = 59.55718233333332
$ws.Range("I8").Value = 0.2958282940622093
$ws.Range("J8").Value = 0.2958282940622093
$ws.Range("M8").Value = 20.56839166666667
$ws.Range("N8").Value = 61.705175
$ws.Range("O8").Value = 0.1304525281245593
$ws.Range("P8").Value = 0.1304525281245593
$ws.Range("Q8").Value = 1224.99545279508
$ws.Range("R8").Value = 11024.95907515572
$ws.Range("S8").Value = 0.03859154885119075
$ws.Range("T8").Value = 0.03859154885119076

# Row 9
$ws.Range("G9").Value = 59.55718233333332
$ws.Range("I9").Value = 0.2958282940622093
$ws.Range("J9").Value = 0.2958282940622093
$ws.Range("O9").Value = 0.6526310778549473
$ws.Range("P9").Value = 0.6526310778549473
$ws.Range("Q9").Value = 6128.437020106722
$ws.Range("R9").Value = 55155.9331809605
$ws.Range("S9").Value = 0.19306673841381
$ws.Range("T9").Value = 0.19306673841381

# Row 10
$ws.Range("G10").Value = 59.55718233333332
$ws.Range("I10").Value = 0.2958282940622093
$ws.Range("J10").Value = 0.2958282940622093
$ws.Range("O10").Value = 0.2169163940204933
$ws.Range("P10").Value = 0.2169163940204934
$ws.Range("S10").Value = 0.06417000679720856
$ws.Range("T10").Value = 0.06417000679720856

